$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I25:I36").NumberFormat = "@"
$ws.Range("Y25:Y36").NumberFormat = "@"
$ws.Range("Z25:Z36").NumberFormat = "@"
$ws.Range("AA25:AA36").NumberFormat = "@"
$ws.Range("AB25:AB36").NumberFormat = "@"

# Row 25
$ws.Range("A25").Value = 112183145
$ws.Range("B25").Value = 96348
$ws.Range("C25").Value = "Ovaliderad"
$ws.Range("D25").Value = "VU"
$ws.Range("E25").Value = 220787
$ws.Range("F25").Value = "Knärot"
$ws.Range("G25").Value = "Goodyera repens"
$ws.Range("H25").Value = "(L.) R. Br."
$ws.Range("I25").Value = "1"
$ws.Range("P25").Value = "Hofors, Gstr"
$ws.Range("Q25").Value = 572350.5295246423
$ws.Range("R25").Value = 6714907.161790377
$ws.Range("S25").Value = 5
$ws.Range("T25").Value = "Gävleborg"
$ws.Range("U25").Value = "Hofors"
$ws.Range("V25").Value = "Gästrikland"
$ws.Range("W25").Value = "Torsåker"
$ws.Range("Y25").Value = "2023-08-23"
$ws.Range("Z25").Value = "00:00"
$ws.Range("AA25").Value = "2023-08-23"
$ws.Range("AB25").Value = "00:00"
$ws.Range("AC25").Value = "Påträffad under Sveaskogs naturvärdesinventering"
$ws.Range("AD25").Value = $false
$ws.Range("AE25").Value = $false
$ws.Range("AG25").Value = $false
$ws.Range("AW25").Value = "Mimmi Persson"
$ws.Range("AX25").Value = "Mimmi Persson"

# Row 26
$ws.Range("A26").Value = 112183148
$ws.Range("B26").Value = 96348
$ws.Range("C26").Value = "Ovaliderad"
$ws.Range("D26").Value = "VU"
$ws.Range("E26").Value = 220787
$ws.Range("F26").Value = "Knärot"
$ws.Range("G26").Value = "Goodyera repens"
$ws.Range("H26").Value = "(L.) R. Br."
$ws.Range("I26").Value = "1"
$ws.Range("P26").Value = "Hofors, Gstr"
$ws.Range("Q26").Value = 572357.010226473
$ws.Range("R26").Value = 6714903.357076311
$ws.Range("S26").Value = 5
$ws.Range("T26").Value = "Gävleborg"
$ws.Range("U26").Value = "Hofors"
$ws.Range("V26").Value = "Gästrikland"
$ws.Range("W26").Value = "Torsåker"
$ws.Range("Y26").Value = "2023-08-23"
$ws.Range("Z26").Value = "00:00"
$ws.Range("AA26").Value = "2023-08-23"
$ws.Range("AB26").Value = "00:00"
$ws.Range("AC26").Value = "Påträffad under Sveaskogs naturvärdesinventering"
$ws.Range("AD26").Value = $false
$ws.Range("AE26").Value = $false
$ws.Range("AG26").Value = $false
$ws.Range("AW26").Value = "Mimmi Persson"
$ws.Range("AX26").Value = "Mimmi Persson"

# Row 27
$ws.Range("A27").Value = 112183149
$ws.Range("B27").Value = 96348
$ws.Range("C27").Value = "Ovaliderad"
$ws.Range("D27").Value = "VU"
$ws.Range("E27").Value = 220787
$ws.Range("F27").Value = "Knärot"
$ws.Range("G27").Value = "Goodyera repens"
$ws.Range("H27").Value = "(L.) R. Br."
$ws.Range("I27").Value = "1"
$ws.Range("P27").Value = "Hofors, Gstr"
$ws.Range("Q27").Value = 572344.933659862
$ws.Range("R27").Value = 6714965.086473988
$ws.Range("S27").Value = 5
$ws.Range("T27").Value = "Gävleborg"
$ws.Range("U27").Value = "Hofors"
$ws.Range("V27").Value = "Gästrikland"
$ws.Range("W27").Value = "Torsåker"
$ws.Range("Y27").Value = "2023-08-23"
$ws.Range("Z27").Value = "00:00"
$ws.Range("AA27").Value = "2023-08-23"
$ws.Range("AB27").Value = "00:00"
$ws.Range("AC27").Value = "Påträffad under Sveaskogs naturvärdesinventering"
$ws.Range("AD27").Value = $false
$ws.Range("AE27").Value = $false
$ws.Range("AG27").Value = $false
$ws.Range("AW27").Value = "Mimmi Persson"
$ws.Range("AX27").Value = "Mimmi Persson"

# Row 28
$ws.Range("A28").Value = 112183137
$ws.Range("B28").Value = 96348
$ws.Range("C28").Value = "Ovaliderad"
$ws.Range("D28").Value = "VU"
$ws.Range("E28").Value = 220787
$ws.Range("F28").Value = "Knärot"
$ws.Range("G28").Value = "Goodyera repens"
$ws.Range("H28").Value = "(L.) R. Br."
$ws.Range("I28").Value = "1"
$ws.Range("P28").Value = "Hofors, Gstr"
$ws.Range("Q28").Value = 572353.8764406883
$ws.Range("R28").Value = 6714961.331215038
$ws.Range("S28").Value = 5
$ws.Range("T28").Value = "Gävleborg"
$ws.Range("U28").Value = "Hofors"
$ws.Range("V28").Value = "Gästrikland"
$ws.Range("W28").Value = "Torsåker"
$ws.Range("Y28").Value = "2023-08-23"
$ws.Range("Z28").Value = "00:00"
$ws.Range("AA28").Value = "2023-08-23"
$ws.Range("AB28").Value = "00:00"
$ws.Range("AC28").Value = "Påträffad under Sveaskogs naturvärdesinventering"
$ws.Range("AD28").Value = $false
$ws.Range("AE28").Value = $false
$ws.Range("AG28").Value = $false
$ws.Range("AW28").Value = "Mimmi Persson"
$ws.Range("AX28").Value = "Mimmi Persson"

# Row 29
$ws.Range("A29").Value = 112183134
$ws.Range("B29").Value = 96348
$ws.Range("C29").Value = "Ovaliderad"
$ws.Range("D29").Value = "VU"
$ws.Range("E29").Value = 220787
$ws.Range("F29").Value = "Knärot"
$ws.Range("G29").Value = "Goodyera repens"
$ws.Range("H29").Value = "(L.) R. Br."
$ws.Range("I29").Value = "1"
$ws.Range("P29").Value = "Hofors, Gstr"
$ws.Range("Q29").Value = 572354.2307759319
$ws.Range("R29").Value = 6714968.224010544
$ws.Range("S29").Value = 5
$ws.Range("T29").Value = "Gävleborg"
$ws.Range("U29").Value = "Hofors"
$ws.Range("V29").Value = "Gästrikland"
$ws.Range("W29").Value = "Torsåker"
$ws.Range("Y29").Value = "2023-08-23"
$ws.Range("Z29").Value = "00:00"
$ws.Range("AA29").Value = "2023-08-23"
$ws.Range("AB29").Value = "00:00"
$ws.Range("AC29").Value = "Påträffad under Sveaskogs naturvärdesinventering"
$ws.Range("AD29").Value = $false
$ws.Range("AE29").Value = $false
$ws.Range("AG29").Value = $false
$ws.Range("AW29").Value = "Mimmi Persson"
$ws.Range("AX29").Value = "Mimmi Persson"

# Row 30
$ws.Range("A30").Value = 112183140
$ws.Range("B30").Value = 96348
$ws.Range("C30").Value = "Ovaliderad"
$ws.Range("D30").Value = "VU"
$ws.Range("E30").Value = 220787
$ws.Range("F30").Value = "Knärot"
$ws.Range("G30").Value = "Goodyera repens"
$ws.Range("H30").Value = "(L.) R. Br."
$ws.Range("I30").Value = "1"
$ws.Range("P30").Value = "Hofors, Gstr"
$ws.Range("Q30").Value = 572349.9172128371
$ws.Range("R30").Value = 6714962.235448033
$ws.Range("S30").Value = 5
$ws.Range("T30").Value = "Gävleborg"
$ws.Range("U30").Value = "Hofors"
$ws.Range("V30").Value = "Gästrikland"
$ws.Range("W30").Value = "Torsåker"
$ws.Range("Y30").Value = "2023-08-23"
$ws.Range("Z30").Value = "00:00"
$ws.Range("AA30").Value = "2023-08-23"
$ws.Range("AB30").Value = "00:00"
$ws.Range("AC30").Value = "Påträffad under Sveaskogs naturvärdesinventering"
$ws.Range("AD30").Value = $false
$ws.Range("AE30").Value = $false
$ws.Range("AG30").Value = $false
$ws.Range("AW30").Value = "Mimmi Persson"
$ws.Range("AX30").Value = "Mimmi Persson"

# Row 31
$ws.Range("A31").Value = 112183141
$ws.Range("B31").Value = 96348
$ws.Range("C31").Value = "Ovaliderad"
$ws.Range("D31").Value = "VU"
$ws.Range("E31").Value = 220787
$ws.Range("F31").Value = "Knärot"
$ws.Range("G31").Value = "Goodyera repens"
$ws.Range("H31").Value = "(L.) R. Br."
$ws.Range("I31").Value = "1"
$ws.Range("P31").Value = "Hofors, Gstr"
$ws.Range("Q31").Value = 572360.8881189874
$ws.Range("R31").Value = 6714980.161596241
$ws.Range("S31").Value = 5
$ws.Range("T31").Value = "Gävleborg"
$ws.Range("U31").Value = "Hofors"
$ws.Range("V31").Value = "Gästrikland"
$ws.Range("W31").Value = "Torsåker"
$ws.Range("Y31").Value = "2023-08-23"
$ws.Range("Z31").Value = "00:00"
$ws.Range("AA31").Value = "2023-08-23"
$ws.Range("AB31").Value = "00:00"
$ws.Range("AC31").Value = "Påträffad under Sveaskogs naturvärdesinventering"
$ws.Range("AD31").Value = $false
$ws.Range("AE31").Value = $false
$ws.Range("AG31").Value = $false
$ws.Range("AW31").Value = "Mimmi Persson"
$ws.Range("AX31").Value = "Mimmi Persson"

# Row 32
$ws.Range("A32").Value = 112183143
$ws.Range("B32").Value = 96348
$ws.Range("C32").Value = "Ovaliderad"
$ws.Range("D32").Value = "VU"
$ws.Range("E32").Value = 220787
$ws.Range("F32").Value = "Knärot"
$ws.Range("G32").Value = "Goodyera repens"
$ws.Range("H32").Value = "(L.) R. Br."
$ws.Range("I32").Value = "1"
$ws.Range("P32").Value = "Hofors, Gstr"
$ws.Range("Q32").Value = 572358.9405480863
$ws.Range("R32").Value = 6714905.363180133
$ws.Range("S32").Value = 5
$ws.Range("T32").Value = "Gävleborg"
$ws.Range("U32").Value = "Hofors"
$ws.Range("V32").Value = "Gästrikland"
$ws.Range("W32").Value = "Torsåker"
$ws.Range("Y32").Value = "2023-08-23"
$ws.Range("Z32").Value = "00:00"
$ws.Range("AA32").Value = "2023-08-23"
$ws.Range("AB32").Value = "00:00"
$ws.Range("AC32").Value = "Påträffad under Sveaskogs naturvärdesinventering"
$ws.Range("AD32").Value = $false
$ws.Range("AE32").Value = $false
$ws.Range("AG32").Value = $false
$ws.Range("AW32").Value = "Mimmi Persson"
$ws.Range("AX32").Value = "Mimmi Persson"

# Row 33
$ws.Range("A33").Value = 112183150
$ws.Range("B33").Value = 96348
$ws.Range("C33").Value = "Ovaliderad"
$ws.Range("D33").Value = "VU"
$ws.Range("E33").Value = 220787
$ws.Range("F33").Value = "Knärot"
$ws.Range("G33").Value = "Goodyera repens"
$ws.Range("H33").Value = "(L.) R. Br."
$ws.Range("I33").Value = "1"
$ws.Range("P33").Value = "Hofors, Gstr"
$ws.Range("Q33").Value = 572358.0913486973
$ws.Range("R33").Value = 6714972.236145046
$ws.Range("S33").Value = 5
$ws.Range("T33").Value = "Gävleborg"
$ws.Range("U33").Value = "Hofors"
$ws.Range("V33").Value = "Gästrikland"
$ws.Range("W33").Value = "Torsåker"
$ws.Range("Y33").Value = "2023-08-23"
$ws.Range("Z33").Value = "00:00"
$ws.Range("AA33").Value = "2023-08-23"
$ws.Range("AB33").Value = "00:00"
$ws.Range("AC33").Value = "Påträffad under Sveaskogs naturvärdesinventering"
$ws.Range("AD33").Value = $false
$ws.Range("AE33").Value = $false
$ws.Range("AG33").Value = $false
$ws.Range("AW33").Value = "Mimmi Persson"
$ws.Range("AX33").Value = "Mimmi Persson"

# Row 34
$ws.Range("A34").Value = 112183151
$ws.Range("B34").Value = 96348
$ws.Range("C34").Value = "Ovaliderad"
$ws.Range("D34").Value = "VU"
$ws.Range("E34").Value = 220787
$ws.Range("F34").Value = "Knärot"
$ws.Range("G34").Value = "Goodyera repens"
$ws.Range("H34").Value = "(L.) R. Br."
$ws.Range("I34").Value = "1"
$ws.Range("P34").Value = "Hofors, Gstr"
$ws.Range("Q34").Value = 572360.9275804338
$ws.Range("R34").Value = 6714978.195056892
$ws.Range("S34").Value = 5
$ws.Range("T34").Value = "Gävleborg"
$ws.Range("U34").Value = "Hofors"
$ws.Range("V34").Value = "Gästrikland"
$ws.Range("W34").Value = "Torsåker"
$ws.Range("Y34").Value = "2023-08-23"
$ws.Range("Z34").Value = "00:00"
$ws.Range("AA34").Value = "2023-08-23"
$ws.Range("AB34").Value = "00:00"
$ws.Range("AC34").Value = "Påträffad under Sveaskogs naturvärdesinventering"
$ws.Range("AD34").Value = $false
$ws.Range("AE34").Value = $false
$ws.Range("AG34").Value = $false
$ws.Range("AW34").Value = "Mimmi Persson"
$ws.Range("AX34").Value = "Mimmi Persson"

# Row 35
$ws.Range("A35").Value = 112183147
$ws.Range("B35").Value = 96348
$ws.Range("C35").Value = "Ovaliderad"
$ws.Range("D35").Value = "VU"
$ws.Range("E35").Value = 220787
$ws.Range("F35").Value = "Knärot"
$ws.Range("G35").Value = "Goodyera repens"
$ws.Range("H35").Value = "(L.) R. Br."
$ws.Range("I35").Value = "1"
$ws.Range("P35").Value = "Hofors, Gstr"
$ws.Range("Q35").Value = 572350.8641464638
$ws.Range("R35").Value = 6714915.037975621
$ws.Range("S35").Value = 5
$ws.Range("T35").Value = "Gävleborg"
$ws.Range("U35").Value = "Hofors"
$ws.Range("V35").Value = "Gästrikland"
$ws.Range("W35").Value = "Torsåker"
$ws.Range("Y35").Value = "2023-08-23"
$ws.Range("Z35").Value = "00:00"
$ws.Range("AA35").Value = "2023-08-23"
$ws.Range("AB35").Value = "00:00"
$ws.Range("AC35").Value = "Påträffad under Sveaskogs naturvärdesinventering"
$ws.Range("AD35").Value = $false
$ws.Range("AE35").Value = $false
$ws.Range("AG35").Value = $false
$ws.Range("AW35").Value = "Mimmi Persson"
$ws.Range("AX35").Value = "Mimmi Persson"

# Row 36
$ws.Range("A36").Value = 112183146
$ws.Range("B36").Value = 96348
$ws.Range("C36").Value = "Ovaliderad"
$ws.Range("D36").Value = "VU"
$ws.Range("E36").Value = 220787
$ws.Range("F36").Value = "Knärot"
$ws.Range("G36").Value = "Goodyera repens"
$ws.Range("H36").Value = "(L.) R. Br."
$ws.Range("I36").Value = "1"
$ws.Range("P36").Value = "Hofors, Gstr"
$ws.Range("Q36").Value = 572346.3926985958
$ws.Range("R36").Value = 6714916.915632093
$ws.Range("S36").Value = 5
$ws.Range("T36").Value = "Gävleborg"
$ws.Range("U36").Value = "Hofors"
$ws.Range("V36").Value = "Gästrikland"
$ws.Range("W36").Value = "Torsåker"
$ws.Range("Y36").Value = "2023-08-23"
$ws.Range("Z36").Value = "00:00"
$ws.Range("AA36").Value = "2023-08-23"
$ws.Range("AB36").Value = "00:00"
$ws.Range("AC36").Value = "Påträffad under Sveaskogs naturvärdesinventering"
$ws.Range("AD36").Value = $false
$ws.Range("AE36").Value = $false
$ws.Range("AG36").Value = $false
$ws.Range("AW36").Value = "Mimmi Persson"
$ws.Range("AX36").Value = "Mimmi Persson"
